$wb = $excel.ActiveWorkbook

# Update the "展览" sheet (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 7006
$ws1.Range("F5").Value = 88

# Update the "全部类型" sheet (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 7006
$ws4.Range("F5").Value = 88
